$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Shelby J" value from B7 up to B6, and clear B7.
$ws.Range("B6").Value = $ws.Range("B7").Value()
$ws.Range("B7").ClearContents()

# Row 6 becomes the tall "author" row (matches the pattern used by row 4, 9, 11),
# row 7 goes back to the default (unset) row height.
$ws.Rows(6).RowHeight = 150
$ws.Rows(7).AutoFit()

# Update the view: scroll so row 5 is the top-left visible row, and move
# the active selection to B13.
$ws.Range("B13").Select()
$excel.ActiveWindow.ScrollRow = 5

$wb.Save()
